# RPA datasets push 2024-07-17
# Insert a new IPO record ("엠83") as row 3, shifting the existing rows down,
# and drop the oldest record ("하스") that falls off the bottom of the table
# so the sheet keeps its original row count (A1:F21).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new row right after "에이치이엠파마(구.에이치이엠)" (row 2).
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row with the new offering's data.
$ws.Range("A3").Value = "엠83"
$ws.Range("B3").Value = "2024.08.01~08.07"
$ws.Range("C3").Value = "11,000~13,000"
$ws.Range("D3").Value = "-"
$ws.Range("E3").Value = 16500
$ws.Range("F3").Value = "신영증권,유진투자증권"

# Drop the row that is now pushed past the end of the table (old "하스" row).
$ws.Rows.Item(22).Delete()
